# Cloudflare DNS export auto-update: a new DNS record
# (check.irrazionale.org) was added to the Cloudflare zone between
# exports, so the export gained one more row. In the sheet this shows
# up as a brand-new row 44 with every following row shifting down by
# one (old row 44 -> new row 45, ..., old row 144 -> new row 145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 44 downward by inserting a fresh row at
# position 44 (mirrors Excel's own "Insert Sheet Rows" behaviour).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row with the new DNS record's data.
$ws.Cells.Item(44, 1).Value  = "2628b62f7c721d0e34d5ad5f81f94aa6"
$ws.Cells.Item(44, 2).Value  = "check.irrazionale.org"
$ws.Cells.Item(44, 3).Value  = "CNAME"
$ws.Cells.Item(44, 4).Value  = "eb879ac1-3924-4f69-9c6c-75f6a9e1cdb8.cfargotunnel.com"
$ws.Cells.Item(44, 5).Value  = $true
$ws.Cells.Item(44, 6).Value  = $true
$ws.Cells.Item(44, 7).Value  = 1
$ws.Cells.Item(44, 8).Value  = "{}"
$ws.Cells.Item(44, 9).Value  = "{}"
$ws.Cells.Item(44, 10).Value = ""
$ws.Cells.Item(44, 11).Value = "[]"
$ws.Cells.Item(44, 12).Value = "2025-04-10T17:26:36.959974Z"
$ws.Cells.Item(44, 13).Value = "2025-04-10T17:26:36.959974Z"
$ws.Cells.Item(44, 14).Value = ""
$ws.Cells.Item(44, 15).Value = ""
